$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: rename first column from "Med" to "Name" ---
$ws.Range("A1").Value = "Name"

# --- Row 2: Atenolol (Ate) / 10 mg Tablet / 1 / 1 / 1
#     Store Code / TGP Code become real numbers instead of text, and
#     "Is Sold" becomes a literal 1 instead of the TRUE() formula ---
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1

# --- Row 3: Colchisin (Colchi) / 500 mg Tablet / 2 / 2 / 1 ---
$ws.Range("A3").Value = "Colchisin (Colchi)"
$ws.Range("B3").Value = "500 mg Tablet"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 1

# --- Row 4: Colchisin (Colchi) / 20 cc Injection / 3 / 3 / 1 ---
$ws.Range("A4").Value = "Colchisin (Colchi)"
$ws.Range("B4").Value = "20 cc Injection"
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 1

# --- New row 5: Atenolol (Ate) / 20 cc Injection / 7 / 7 / 1 ---
$ws.Range("A5").Value = "Atenolol (Ate)"
$ws.Range("B5").Value = "20 cc Injection"
$ws.Range("C5").Value = 7
$ws.Range("D5").Value = 7
$ws.Range("E5").Value = 1

# --- New row 6: Losartan (Lora) / 10 mg Tablet / 5 / 5 / 1 ---
$ws.Range("A6").Value = "Losartan (Lora)"
$ws.Range("B6").Value = "10 mg Tablet"
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 5
$ws.Range("E6").Value = 1

# --- Make the formatting of the (now reordered/renumbered) data rows
#     consistent: rows 2-4 use wrapped text like the rest of that block ---
$ws.Range("A3").WrapText = $true
$ws.Range("A4:E4").WrapText = $true

# --- Match the cursor/selection position left behind after entering the
#     new rows ---
$ws.Range("B12").Select() | Out-Null
